$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[Wei-jie%Guan%NULL%0,               Zheng-yi%Ni%NULL%0,               Zheng-yi%Ni%NULL%0,               Yu%Hu%NULL%0,               Wen-hua%Liang%NULL%0,               Chun-quan%Ou%NULL%0,               Jian-xing%He%NULL%0,               Lei%Liu%NULL%0,               Hong%Shan%NULL%0,               Chun-liang%Lei%NULL%0,               David S.C.%Hui%NULL%0,               Bin%Du%NULL%0,               Lan-juan%Li%NULL%0,               Guang%Zeng%NULL%0,               Kwok-Yung%Yuen%NULL%0,               Ru-chong%Chen%NULL%0,               Chun-li%Tang%NULL%0,               Tao%Wang%NULL%0,               Ping-yan%Chen%NULL%0,               Jie%Xiang%NULL%0,               Shi-yue%Li%NULL%0,               Jin-lin%Wang%NULL%0,               Zi-jing%Liang%NULL%0,               Yi-xiang%Peng%NULL%0,               Li%Wei%NULL%0,               Yong%Liu%NULL%0,               Ya-hua%Hu%NULL%0,               Peng%Peng%NULL%0,               Jian-ming%Wang%NULL%0,               Ji-yang%Liu%NULL%0,               Zhong%Chen%NULL%0,               Gang%Li%NULL%0,               Zhi-jian%Zheng%NULL%0,               Shao-qin%Qiu%NULL%0,               Jie%Luo%NULL%0,               Chang-jiang%Ye%NULL%0,               Shao-yong%Zhu%NULL%0,               Nan-shan%Zhong%NULL%0]'
$ws.Range("E3").Value2 = '[Ping%Wu%NULL%1,               Fang%Duan%NULL%1,               Chunhua%Luo%NULL%1,               Qiang%Liu%NULL%0,               Xingguang%Qu%NULL%1,               Liang%Liang%NULL%1,               Kaili%Wu%NULL%1]'
$ws.Range("E4").Value2 = '[Xian%Zhang%NULL%1,               Xuhui%Chen%NULL%3,               Liwen%Chen%NULL%1,               Chaohua%Deng%NULL%3,               Xiaojing%Zou%NULL%0,               Xiaojing%Zou%NULL%0,               Weiyong%Liu%NULL%3,               Weiyong%Liu%NULL%0,               Huimin%Yu%NULL%3,               Huimin%Yu%NULL%0,               Bo%Chen%NULL%2,               Xufang%Sun%NULL%2]'
$ws.Range("E5").Value2 = '[Jianhua%Xia%NULL%1,               Jianping%Tong%NULL%1,               Mengyun%Liu%NULL%1,               Ye%Shen%idrshen@zju.edu.cn%1,               Dongyu%Guo%guodongyu666@sina.com%1]'
$ws.Range("E6").Value2 = '[Hua-Tao%Xie%NULL%1,               Shi-Yun%Jiang%NULL%1,               Kang-Kang%Xu%NULL%1,               Xin%Liu%NULL%1,               Bing%Xu%NULL%1,               Lin%Wang%lin_wang@hust.edu.cn%0,               Ming-Chang%Zhang%mingchangzhang@hotmail.com%1]'
$ws.Range("E7").Value2 = '[Liwen%Chen%NULL%1,               Chaohua%Deng%NULL%0,               Chaohua%Deng%NULL%0,               Xuhui%Chen%NULL%0,               Xuhui%Chen%NULL%0,               Xian%Zhang%NULL%1,               Bo%Chen%NULL%0,               Huimin%Yu%NULL%0,               Yuanjun%Qin%NULL%1,               Ke%Xiao%NULL%1,               Hong%Zhang%tjyksys@163.com%1,               Xufang%Sun%sunxufang2016@163.com%0]'
$ws.Range("E8").Value2 = '[Kiran%Kumar%NULL%1,               Akshata A%Prakash%NULL%1,               Suresh Babu%Gangasagara%NULL%1,               Sujatha BL%Rathod%NULL%1,               K%Ravi%NULL%1,               Ambica%Rangaiah%NULL%1,               Sathyanarayan Muthur%Shankar%NULL%1,               Shantala Gowdara%Basawarajappa%NULL%1,               Shashi%Bhushan%NULL%1,               T G%Neeraja%NULL%1,               Srinivas%Khandenahalli%NULL%1,               M%Swetha%NULL%1,               Priyam%Gupta%NULL%1,               U C%Sampritha%NULL%1,               Guru N S%Prasad%NULL%1,               Chakravarthy Raghunathan%Jayanthi%NULL%1]'
$ws.Range("E9").Value2 = '[Yunyun%Zhou%NULL%1,               Chaoye%Duan%NULL%1,               Yuyang%Zeng%NULL%1,               Yongqing%Tong%NULL%1,               Yuhong%Nie%NULL%1,               Yang%Yang%NULL%0,               Zhen%Chen%NULL%1,               Changzheng%Chen%chenchangzheng@whu.edu.cn%1]'
$ws.Range("E10").Value2 = '[Saeed%Karimi%NULL%1,               Amir%Arabi%amir_arab_91@yahoo.com%1,               Toktam%Shahraki%NULL%2,               Toktam%Shahraki%NULL%0,               Sare%Safi%NULL%1]'
$ws.Range("E11").Value2 = '[Marion%Grimaud%NULL%2,               Julie%Starck%NULL%1,               Michael%Levy%NULL%1,               Clémence%Marais%NULL%1,               Judith%Chareyre%NULL%1,               Diala%Khraiche%NULL%1,               Marianne%Leruez-Ville%NULL%1,               Pierre%Quartier%NULL%1,               Pierre Louis%Léger%NULL%1,               Guillaume%Geslain%NULL%1,               Nada%Semaan%NULL%1,               Florence%Moulin%NULL%1,               Matthieu%Bendavid%NULL%1,               Sandrine%Jean%NULL%1,               Géraldine%Poncelet%NULL%1,               Sylvain%Renolleau%NULL%1,               Mehdi%Oualha%mehdi.oualha@aphp.fr%1]'
$ws.Range("E12").Value2 = '[Zhixiong%Fang%NULL%1,               Yi%Zhang%NULL%0,               Changfa%Hang%NULL%1,               Jingwen%Ai%NULL%1,               Shaojie%Li%NULL%1,               Wenhong%Zhang%NULL%2]'
$ws.Range("E13").Value2 = '[Paola%Valente%NULL%1,               Giancarlo%Iarossi%NULL%1,               Matteo%Federici%NULL%1,               Sergio%Petroni%NULL%1,               Paolo%Palma%NULL%1,               Nicola%Cotugno%NULL%1,               Maria A.%De Ioris%NULL%1,               Andrea%Campana%NULL%1,               Luca%Buzzonetti%NULL%1]'
$ws.Range("E14").Value2 = '[Ivan Yu Jun%Seah%NULL%1,               Danielle E.%Anderson%NULL%1,               Adrian Eng Zheng%Kang%NULL%1,               Linfa%Wang%NULL%1,               Pooja%Rao%NULL%1,               Barnaby Edward%Young%NULL%0,               David Chien%Lye%NULL%1,               Rupesh%Agrawal%Rupesh_agrawal@ttsh.com.sg%1]'
